$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-12: 45233 -> 45243
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
